# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values replacing the old "Strike#" values in column G (rows 2-20)
$kValues = @{
    2  = 2
    3  = 6
    4  = 3
    5  = 2
    6  = 7
    7  = 3
    8  = 6
    9  = 8
    10 = 1
    11 = 3
    12 = 9
    13 = 5
    14 = 6
    15 = 3
    16 = 5
    17 = 12
    18 = 5
    19 = 1
    20 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
